# Finished phylo signal analyses: refresh the summary table (columns J:P)
# with re-run values and give the new "% of phylogenies with significant
# signal" column (P) a one-decimal number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated summary statistics (columns N, O, P for rows 3-7) ---------

# Row 3
$ws.Range("N3").Value = 1.62941832060485
$ws.Range("O3").Value = 0.308198154861495
$ws.Range("P3").Value = 4.4943820224719104

# Row 4
$ws.Range("N4").Value = 0.918518977061572
$ws.Range("O4").Value = 0.332061735862865
$ws.Range("P4").Value = 8.9887640449438209

# Row 5
$ws.Range("N5").Value = 0.0803091516668518
$ws.Range("O5").Value = 0.106385467691284
$ws.Range("P5").Value = 0

# Row 6
$ws.Range("N6").Value = 0.696652403902864
$ws.Range("O6").Value = 0.0382005693309071
$ws.Range("P6").Value = 100

# Row 7
$ws.Range("N7").Value = 0.377407239940525
$ws.Range("O7").Value = 0.0951302234275259
$ws.Range("P7").Value = 82.022471910112401

# Column P now carries its own one-decimal-place number format (a new
# style distinct from N/O's existing "0.000" style).
$ws.Range("P3:P7").NumberFormat = "0.0"

# --- Column width tweaks -------------------------------------------------
# Column K (11) was manually narrowed (loses best-fit autosizing).
$ws.Columns.Item(11).ColumnWidth = 14.833333333333332

# Columns N (14) and O (15) used to share one best-fit column-width entry;
# they are now split into two separate best-fit columns with narrower
# widths reflecting the shorter re-computed values.
$ws.Columns.Item(14).ColumnWidth = 7.0
$ws.Columns.Item(15).ColumnWidth = 7.500000000000001

# --- Selection ------------------------------------------------------------
$ws.Range("J3:P7").Select()
